# Automatische test-sync: 2025-06-20 09:30:50
# Adds the new "Sollicitatie marketingfunctie" mail-log row (row 5) to the
# "Logs" sheet, the matching dashboard aggregate row, extends the
# conditional-formatting ranges to include the new row, and widens the
# chart's category/value series references accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Logs" sheet - append the new mail entry as row 5
# ---------------------------------------------------------------------
$wsLogs = $wb.Worksheets.Item("Logs")

$wsLogs.Range("A5").Value = "Sollicitatie marketingfunctie"
$wsLogs.Range("B5").Value = "mailmind.test@zohomail.eu"
$wsLogs.Range("C5").Value = "Hierbij solliciteer ik voor de functie van marketeer. Zie bijlage voor CV."
$wsLogs.Range("D5").Value = "Sollicitatie / Vacature"
$wsLogs.Range("F5").Value = "2025-06-20 09:30:14"
$wsLogs.Range("G5").Value = "Nee"

# Extend the existing conditional-formatting rules (category + answered
# columns) so they keep covering rows 2-5 instead of 2-4. Using
# ModifyAppliesToRange preserves each rule's dxfId/priority/formula -
# only the applicable range (sqref) changes.
foreach ($cf in $wsLogs.Range("D2:D4").FormatConditions) {
    $cf.ModifyAppliesToRange($wsLogs.Range("D2:D5"))
}
foreach ($cf in $wsLogs.Range("G2:G4").FormatConditions) {
    $cf.ModifyAppliesToRange($wsLogs.Range("G2:G5"))
}

# ---------------------------------------------------------------------
# 2. "Dashboard" sheet - append the matching category total as row 5
# ---------------------------------------------------------------------
$wsDash = $wb.Worksheets.Item("Dashboard")

$wsDash.Range("A5").Value = "Sollicitatie / Vacature"
$wsDash.Range("B5").Value = 1

# ---------------------------------------------------------------------
# 3. Chart - widen the category/value series ranges from row 4 to row 5
# ---------------------------------------------------------------------
$chartObj = $wsDash.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES('Dashboard'!B1,'Dashboard'!`$A`$2:`$A`$5,'Dashboard'!`$B`$2:`$B`$5,1)"
